# Generate Report for Handback
# Mark the d86bf385-...-.md file (row 3 in every sheet) as handed back:
# update its Status/locale columns, set the new handback timestamp, and
# clear the stale "handback not latest" error message.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $handedBack
$wsZhCn.Range("K3").Value = "2016-09-03 16:55:31"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $handedBack
$wsDeDe.Range("K3").Value = "2016-09-03 16:55:38"
$wsDeDe.Range("P3").Value = ""
